$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Catégorie" (SEO) column for rows 11 and 12, which were
# previously left blank in column A.
$ws.Range("A11").Value = "SEO"
$ws.Range("A12").Value = "SEO"

# Add the new "Accessibilité" row (row 13) with its four associated
# descriptive cells (problem, explanation, best practice, recommended
# action).
$ws.Range("A13").Value = "Accessibilité"
$ws.Range("B13").Value = "mauvais ratio de couleur"
$ws.Range("C13").Value = "difficulté de visibilité"
$ws.Range("D13").Value = "adopter les bons ratio"
$ws.Range("E13").Value = "corriger les ratios de couleurs"

# Move the active selection from E12 to B12.
$ws.Range("B12").Select()
